$d = $word.ActiveDocument

# Locate the paragraph that ends the "My" narrative - the one discussing
# "garanterade timmar." - so the two new paragraphs are appended right
# after it (and still before the trailing bookmarkEnd), regardless of
# exact paragraph index.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*garanterade timmar.*") {
        $anchor = $p
    }
}
if ($anchor -eq $null) {
    $anchor = $d.Paragraphs.Last
}

# First new paragraph: statistics from SOU mentioned by Arpi (2012).
# InsertParagraphAfter() on the anchor's range clones the anchor's own
# paragraph formatting (style "BodyText"), so no explicit style push is
# needed - that also avoids stamping stray rsid attributes on the <w:p>.
$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p1 = $anchor.Next()
$p1.Range.Text = "Arpi (2012) inkluderar också statistik ifrån SOU vilka pekar på att bemanningsbranschen tas upp huvudligen av personer där arbetet passar livssituationen. Dessutom så överrepresenteras vissa grupper, som unga, kvinnor, m.m."

# Second new paragraph: closing statement about un-guaranteed hours.
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Slutligen så påstår Arpi (2012) att arbetare på bemanningsföretag borde vara förberedda på ogaranterade arbetstimmar."

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
